$d = $word.ActiveDocument

$pairs = @(
    @("804÷3=268, 0", "926÷6=154, 2"),
    @("232÷4=58, 0", "277÷8=34, 5"),
    @("473÷9=52, 5", "899÷2=449, 1"),
    @("482÷8=60, 2", "549÷5=109, 4"),
    @("400÷5=80, 0", "597÷2=298, 1"),
    @("225÷8=28, 1", "657÷2=328, 1"),
    @("974÷4=243, 2", "990÷9=110, 0"),
    @("612÷2=306, 0", "579÷2=289, 1"),
    @("562÷3=187, 1", "640÷6=106, 4"),
    @("900÷5=180, 0", "668÷9=74, 2"),
    @("273÷7=39, 0", "594÷7=84, 6"),
    @("485÷9=53, 8", "702÷8=87, 6"),
    @("956÷8=119, 4", "856÷8=107, 0"),
    @("581÷5=116, 1", "655÷3=218, 1"),
    @("403÷8=50, 3", "750÷9=83, 3"),
    @("995÷7=142, 1", "390÷6=65, 0"),
    @("968÷9=107, 5", "799÷7=114, 1"),
    @("921÷7=131, 4", "708÷3=236, 0"),
    @("457÷8=57, 1", "569÷9=63, 2"),
    @("820÷2=410, 0", "221÷9=24, 5"),
    @("420÷6=70, 0", "330÷2=165, 0"),
    @("214÷9=23, 7", "110÷7=15, 5"),
    @("965÷9=107, 2", "747÷2=373, 1"),
    @("679÷2=339, 1", "318÷7=45, 3"),
    @("499÷7=71, 2", "624÷5=124, 4")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
